# Applies the "crusader" review edits:
#  1. Inserts a new "Meta description: ..." paragraph right after the
#     top H1 title.
#  2. Removes the duplicate bold title paragraph that had been placed
#     near the end of the document (right before the final italic blurb).
#  3. Replaces the text of that trailing italic paragraph with the new
#     AI feature-image prompt text.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Part 1: add the "Meta description" paragraph right after the title ---
$titlePara = $d.Paragraphs.Item(1)
Write-Host "Title paragraph:" $titlePara.Range.Text

$titlePara.Range.InsertParagraphAfter() | Out-Null

$metaPara = $d.Paragraphs.Item(2)
$metaXml = '<w:p ' + $wNs + '><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Crusader, a medieval-themed online slot game. Play for free and trigger free spins with the wild symbol.</w:t></w:r></w:p>'
$metaPara.Range.InsertXML($metaXml)
Write-Host "Inserted meta paragraph:" $d.Paragraphs.Item(2).Range.Text

# --- Part 2: drop the duplicated "Play Crusader Free..." title near the end,
#     then rewrite the final italic paragraph with the new prompt copy. ---
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($count - 1)
Write-Host "Paragraph to remove (duplicate title):" $dupTitlePara.Range.Text
$dupTitlePara.Range.Delete() | Out-Null

$count = $d.Paragraphs.Count
$blurbPara = $d.Paragraphs.Item($count)
Write-Host "Paragraph to rewrite (old blurb):" $blurbPara.Range.Text

$blurbText = "Create an eye-catching feature image for Crusader that incorporates the game's Medieval theme and features a happy Maya warrior with glasses. Use bright colors to make the image pop and make sure to highlight the warrior's glasses to add a touch of uniqueness. The cartoon style of the image should be playful and inviting, with a dynamic pose for the warrior that exudes confidence. The background should feature symbols that represent medieval warfare, such as swords, shields, and castles. Make sure that the overall design of the image is consistent with the game's theme, while also being fun and engaging."
$blurbXml = '<w:p ' + $wNs + '><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>' + $blurbText + '</w:t></w:r></w:p>'
$blurbPara.Range.InsertXML($blurbXml)

Write-Host "Final paragraph count:" $d.Paragraphs.Count
Write-Host "Final paragraph text:" $d.Paragraphs.Item($d.Paragraphs.Count).Range.Text
